# Jenkins integration pass: rename the "Verifyloginout" test-case label to
# "Verifyloginoutwithalltests" everywhere it appears on the sheet, and tidy
# up the sheet view (columns widened to fit the longer text, selection
# reset to A1) the way Excel leaves the file after the author re-saved it
# post-edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# A1 and D4 both held the old "Verifyloginout" label - update both to the
# new, longer test name used after the Jenkins integration.
$ws.Range("A1").Value = "Verifyloginoutwithalltests"
$ws.Range("D4").Value = "Verifyloginoutwithalltests"

# Columns A and D were widened so the longer label fits (31 and 13 for B/C
# are untouched). ColumnWidth is quantized internally to 1/6-character
# steps by this host, so feed it values that round-trip to the closest
# representable width to the real target (25.21875 / 25.6640625 chars).
$ws.Columns.Item(1).ColumnWidth = 24.33
$ws.Columns.Item(4).ColumnWidth = 24.83

# The old sheet had B1:C1 selected; the saved workbook just keeps the
# default A1 selection, so select A1 explicitly.
$ws.Range("A1").Select()
